# Applies the "Revu et implémentation des specifications" edit:
#  - Renames sheet 3 ("Feuil3" -> "Courbe d'apparition ennemis") and
#    populates its enemy-spawn-curve table/formulas.
#  - Scoring sheet: swaps the LIFE label for POINTS DE VIE, tweaks the
#    enemy-count sample values, adds the extra parentheses to the score
#    formulas and documents them in column C.
#  - Feedback sheet: fixes a typo ("skill" -> "kill") in one message.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Sheet 3: "Feuil3" -> "Courbe d'apparition ennemis"
# ---------------------------------------------------------------------
$ws3.Name = "Courbe d'apparition ennemis"

$ws3.Range("A2").Value = "SCORE"
$ws3.Range("B2").Value = 1178

$ws3.Range("A3").Value = "COEFFICIENT DE DIFFICULTE"
$ws3.Range("B3").Value = 1

$ws3.Range("A4").Value = "POINTS DE VIE"
$ws3.Range("B4").Value = 3

$ws3.Range("A5").Value = "DUREE DE LA PARTIE"
$ws3.Range("B5").Value = 60

$ws3.Range("A9").Value = "Nombre de créatures maximum dans l'arène"
$ws3.Range("B9").Formula = "=((B2+B3)/B5)/2"

$ws3.Range("A10").Value = "Intervalle d'apparition des créatures"
$ws3.Range("B10").Formula = "=B7- (((B9-B6)+B3)/(10-B4))"

$ws3.Range("B9").NumberFormat = "0"

$ws3.Range("A6").Value = "NOMBRES DE CREATURES DANS L'ARENE"
$ws3.Range("B6").Value = 1

$ws3.Range("A7").Value = "INTERVALLE INITIAL"
$ws3.Range("B7").Value = 3

$ws3.Columns.Item(1).ColumnWidth = 40

# ---------------------------------------------------------------------
# Sheet 2: "Feedback"
# ---------------------------------------------------------------------
$ws2.Range("C2").Value = "Score gagner pour le kill de l'ennemi"

# ---------------------------------------------------------------------
# Sheet 1: "Scoring"
# ---------------------------------------------------------------------
$ws1.Range("F2").Value = "POINTS DE VIE"
$ws1.Range("C3").Value = 0
$ws1.Range("C4").Value = 10

$ws1.Range("B6").Formula = "=B3*((AVERAGE(C3:C4)/G2))"
$ws1.Range("C6").Value = "B3*((MOYENNE(C3:C4)/G2)*B5)"

$ws1.Range("B7").Formula = "=B4*((AVERAGE(C3:C4)/G2))"
$ws1.Range("C7").Value = "B4*((MOYENNE(C3:C4)/G2)*B5)"

$ws1.Columns.Item(6).ColumnWidth = 12.74

# ---------------------------------------------------------------------
# Selections — recorded per-sheet, last-activated sheet stays the
# workbook's active/tab-selected sheet (sheet 3, matching "before").
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("C14").Select()

$ws2.Activate()
$ws2.Range("B16").Select()

$ws3.Activate()
$ws3.Range("B9").Select()
